$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 24.04021751172397
$ws.Cells.Item(2, 3).Value = 6.086104322368083
$ws.Cells.Item(2, 4).Value = 6.09197249722454
$ws.Cells.Item(2, 5).Value = 9.372883342665713
$ws.Cells.Item(2, 7).Value = 3.795915713121647
$ws.Cells.Item(2, 11).Value = 20.08986664291081
$ws.Cells.Item(2, 12).Value = 10.75746748680636
$ws.Cells.Item(2, 14).Value = 27.83292940279413
$ws.Cells.Item(3, 2).Value = 23.85119707969623
$ws.Cells.Item(3, 3).Value = 5.886702567013547
$ws.Cells.Item(3, 4).Value = 5.991341331905851
$ws.Cells.Item(3, 5).Value = 9.386260406253207
$ws.Cells.Item(3, 7).Value = 3.800186026044584
$ws.Cells.Item(3, 11).Value = 19.95689306740989
$ws.Cells.Item(3, 12).Value = 10.75416073945871
$ws.Cells.Item(3, 14).Value = 27.74130763278652
$ws.Cells.Item(4, 2).Value = 23.74111252507031
$ws.Cells.Item(4, 3).Value = 5.762968095838743
$ws.Cells.Item(4, 4).Value = 5.930679434801709
$ws.Cells.Item(4, 5).Value = 9.395849270376406
$ws.Cells.Item(4, 7).Value = 3.802941205902651
$ws.Cells.Item(4, 11).Value = 19.88025302162389
$ws.Cells.Item(4, 12).Value = 10.7542499524299
$ws.Cells.Item(4, 14).Value = 27.68588775791095
$ws.Cells.Item(5, 2).Value = 23.69779402995292
$ws.Cells.Item(5, 3).Value = 5.712307193898692
$ws.Cells.Item(5, 4).Value = 5.906275821393876
$ws.Cells.Item(5, 5).Value = 9.40010258033727
$ws.Cells.Item(5, 7).Value = 3.804097593294629
$ws.Cells.Item(5, 11).Value = 19.85030512080923
$ws.Cells.Item(5, 12).Value = 10.754819409846
$ws.Cells.Item(5, 14).Value = 27.66351879536738
$ws.Cells.Item(6, 2).Value = 23.69069524387651
$ws.Cells.Item(6, 3).Value = 5.703883326527391
$ws.Cells.Item(6, 4).Value = 5.902243709024802
$ws.Cells.Item(6, 5).Value = 9.400829720390695
$ws.Cells.Item(6, 7).Value = 3.804291645693322
$ws.Cells.Item(6, 11).Value = 19.84541049307359
$ws.Cells.Item(6, 12).Value = 10.7549461642131
$ws.Cells.Item(6, 14).Value = 27.65981757183428
$ws.Cells.Item(7, 2).Value = 23.74052202402721
$ws.Cells.Item(7, 3).Value = 5.762285707215841
$ws.Cells.Item(7, 4).Value = 5.93034899338044
$ws.Cells.Item(7, 5).Value = 9.395905232093025
$ws.Cells.Item(7, 7).Value = 3.8029566650047
$ws.Cells.Item(7, 11).Value = 19.87984390594673
$ws.Cells.Item(7, 12).Value = 10.75425547386577
$ws.Cells.Item(7, 14).Value = 27.68558520538879
$ws.Cells.Item(8, 2).Value = 23.97382778907017
$ws.Cells.Item(8, 3).Value = 6.017672168633765
$ws.Cells.Item(8, 4).Value = 6.057061364677553
$ws.Cells.Item(8, 5).Value = 9.377210286974771
$ws.Cells.Item(8, 7).Value = 3.79736055626106
$ws.Cells.Item(8, 11).Value = 20.04299379429827
$ws.Cells.Item(8, 12).Value = 10.75588754937755
$ws.Cells.Item(8, 14).Value = 27.80116759240299
$ws.Cells.Item(9, 2).Value = 24.47687421032802
$ws.Cells.Item(9, 3).Value = 6.504565516123921
$ws.Cells.Item(9, 4).Value = 6.312993237126757
$ws.Cells.Item(9, 5).Value = 9.351463900950655
$ws.Cells.Item(9, 7).Value = 3.787437136916314
$ws.Cells.Item(9, 11).Value = 20.4014638925118
$ws.Cells.Item(9, 12).Value = 10.77589124708599
$ws.Cells.Item(9, 14).Value = 28.0342446665679
$ws.Cells.Item(10, 2).Value = 24.87153778607586
$ws.Cells.Item(10, 3).Value = 6.849451466362471
$ws.Cells.Item(10, 4).Value = 6.503582468660117
$ws.Cells.Item(10, 5).Value = 9.339203523290429
$ws.Cells.Item(10, 7).Value = 3.780778044091863
$ws.Cells.Item(10, 11).Value = 20.68659347367799
$ws.Cells.Item(10, 12).Value = 10.8007982245894
$ws.Cells.Item(10, 14).Value = 28.20910699939145
$ws.Cells.Item(11, 2).Value = 25.05591062061714
$ws.Cells.Item(11, 3).Value = 7.002777178105789
$ws.Cells.Item(11, 4).Value = 6.590429777120375
$ws.Cells.Item(11, 5).Value = 9.335070725027325
$ws.Cells.Item(11, 7).Value = 3.777883922859787
$ws.Cells.Item(11, 11).Value = 20.82062690792077
$ws.Cells.Item(11, 12).Value = 10.81433340100223
$ws.Cells.Item(11, 14).Value = 28.2893755555659
$ws.Cells.Item(12, 2).Value = 25.12636960691606
$ws.Cells.Item(12, 3).Value = 7.060268523831398
$ws.Cells.Item(12, 4).Value = 6.62330461239726
$ws.Cells.Item(12, 5).Value = 9.33371332466788
$ws.Cells.Item(12, 7).Value = 3.776807281045161
$ws.Cells.Item(12, 11).Value = 20.87196783139517
$ws.Cells.Item(12, 12).Value = 10.81977427831542
$ws.Cells.Item(12, 14).Value = 28.3198688070378
$ws.Cells.Item(13, 2).Value = 25.11116737792608
$ws.Cells.Item(13, 3).Value = 7.047912918806504
$ws.Cells.Item(13, 4).Value = 6.61622549911221
$ws.Cells.Item(13, 5).Value = 9.333996434438538
$ws.Cells.Item(13, 7).Value = 3.777038298905063
$ws.Cells.Item(13, 11).Value = 20.86088520771176
$ws.Cells.Item(13, 12).Value = 10.81858848995113
$ws.Cells.Item(13, 14).Value = 28.31329732077185
$ws.Cells.Item(14, 2).Value = 25.06169479468186
$ws.Cells.Item(14, 3).Value = 7.007518750119426
$ws.Cells.Item(14, 4).Value = 6.59313484505655
$ws.Cells.Item(14, 5).Value = 9.334954891108218
$ws.Cells.Item(14, 7).Value = 3.777794960868425
$ws.Cells.Item(14, 11).Value = 20.82483923950575
$ws.Cells.Item(14, 12).Value = 10.81477471467247
$ws.Cells.Item(14, 14).Value = 28.29188234270296
$ws.Cells.Item(15, 2).Value = 25.03147324222346
$ws.Cells.Item(15, 3).Value = 6.982700364037439
$ws.Cells.Item(15, 4).Value = 6.578988568665705
$ws.Cells.Item(15, 5).Value = 9.335569004877934
$ws.Cells.Item(15, 7).Value = 3.778260947537682
$ws.Cells.Item(15, 11).Value = 20.80283516503252
$ws.Cells.Item(15, 12).Value = 10.81247968752534
$ws.Cells.Item(15, 14).Value = 28.27877751949839
$ws.Cells.Item(16, 2).Value = 24.859581292905
$ws.Cells.Item(16, 3).Value = 6.839354508597395
$ws.Cells.Item(16, 4).Value = 6.497906863443842
$ws.Cells.Item(16, 5).Value = 9.339502669637609
$ws.Cells.Item(16, 7).Value = 3.780969889580595
$ws.Cells.Item(16, 11).Value = 20.67791801688822
$ws.Cells.Item(16, 12).Value = 10.79995790134358
$ws.Cells.Item(16, 14).Value = 28.20387524191915
$ws.Cells.Item(17, 2).Value = 24.75533178308846
$ws.Cells.Item(17, 3).Value = 6.750460779846724
$ws.Cells.Item(17, 4).Value = 6.448179367331522
$ws.Cells.Item(17, 5).Value = 9.342285750712714
$ws.Cells.Item(17, 7).Value = 3.782666254026981
$ws.Cells.Item(17, 11).Value = 20.60236773092419
$ws.Cells.Item(17, 12).Value = 10.79283982145036
$ws.Cells.Item(17, 14).Value = 28.15810558878685
$ws.Cells.Item(18, 2).Value = 24.6958289957698
$ws.Cells.Item(18, 3).Value = 6.698998944301001
$ws.Cells.Item(18, 4).Value = 6.419593231953595
$ws.Cells.Item(18, 5).Value = 9.344022479578152
$ws.Cells.Item(18, 7).Value = 3.783654685327981
$ws.Cells.Item(18, 11).Value = 20.55932272366294
$ws.Cells.Item(18, 12).Value = 10.78895334008338
$ws.Cells.Item(18, 14).Value = 28.13184854070143
$ws.Cells.Item(19, 2).Value = 24.6757628190391
$ws.Cells.Item(19, 3).Value = 6.68151958754392
$ws.Cells.Item(19, 4).Value = 6.409918213109977
$ws.Cells.Item(19, 5).Value = 9.344633862613344
$ws.Cells.Item(19, 7).Value = 3.783991541168562
$ws.Cells.Item(19, 11).Value = 20.54481986071131
$ws.Cells.Item(19, 12).Value = 10.78767315417799
$ws.Cells.Item(19, 14).Value = 28.12297030989049
$ws.Cells.Item(20, 2).Value = 24.76638223327749
$ws.Cells.Item(20, 3).Value = 6.759958517180578
$ws.Cells.Item(20, 4).Value = 6.453471555099477
$ws.Cells.Item(20, 5).Value = 9.341975414864201
$ws.Cells.Item(20, 7).Value = 3.782484356869241
$ws.Cells.Item(20, 11).Value = 20.61036806701685
$ws.Cells.Item(20, 12).Value = 10.7935760733385
$ws.Cells.Item(20, 14).Value = 28.16297081007502
$ws.Cells.Item(21, 2).Value = 25.07620913674596
$ws.Cells.Item(21, 3).Value = 7.019399384819441
$ws.Cells.Item(21, 4).Value = 6.599917726534627
$ws.Cells.Item(21, 5).Value = 9.334667736116511
$ws.Cells.Item(21, 7).Value = 3.777572188194525
$ws.Cells.Item(21, 11).Value = 20.83541122595104
$ws.Cells.Item(21, 12).Value = 10.81588636767348
$ws.Cells.Item(21, 14).Value = 28.29816985444545
$ws.Cells.Item(22, 2).Value = 25.28240908468831
$ws.Cells.Item(22, 3).Value = 7.185615732990663
$ws.Cells.Item(22, 4).Value = 6.695542076987017
$ws.Cells.Item(22, 5).Value = 9.331101711117455
$ws.Cells.Item(22, 7).Value = 3.774474230317298
$ws.Cells.Item(22, 11).Value = 20.98588413959172
$ws.Cells.Item(22, 12).Value = 10.83230485796762
$ws.Cells.Item(22, 14).Value = 28.38709410111514
$ws.Cells.Item(23, 2).Value = 25.17203515387403
$ws.Cells.Item(23, 3).Value = 7.097226137901886
$ws.Cells.Item(23, 4).Value = 6.644524254556015
$ws.Cells.Item(23, 5).Value = 9.332894305681215
$ws.Cells.Item(23, 7).Value = 3.776117425399597
$ws.Cells.Item(23, 11).Value = 20.90527576953577
$ws.Cells.Item(23, 12).Value = 10.82337449687699
$ws.Cells.Item(23, 14).Value = 28.33958421522008
$ws.Cells.Item(24, 2).Value = 24.76138497613033
$ws.Cells.Item(24, 3).Value = 6.755665693485883
$ws.Cells.Item(24, 4).Value = 6.451078945529725
$ws.Cells.Item(24, 5).Value = 9.342115291907502
$ws.Cells.Item(24, 7).Value = 3.782566551546402
$ws.Cells.Item(24, 11).Value = 20.60674989797088
$ws.Cells.Item(24, 12).Value = 10.79324257272384
$ws.Cells.Item(24, 14).Value = 28.16077106693471
$ws.Cells.Item(25, 2).Value = 24.3361801629801
$ws.Cells.Item(25, 3).Value = 6.374813223184121
$ws.Cells.Item(25, 4).Value = 6.243165706148387
$ws.Cells.Item(25, 5).Value = 9.357260158617489
$ws.Cells.Item(25, 7).Value = 3.790010127781413
$ws.Cells.Item(25, 11).Value = 20.30053414650594
$ws.Cells.Item(25, 12).Value = 10.76868315925992
$ws.Cells.Item(25, 14).Value = 27.97053306643549
